{"js": "// Update the worksheet date and all 25 multiplication problems to the\n// next day's values, per the commit diff.\nconst replacements = [\n  [\"2025-11-15 Saturday\", \"2025-11-16 Sunday\"],\n  [\"810\u00d77=\", \"718\u00d79=\"],\n  [\"659\u00d79=\", \"633\u00d77=\"],\n  [\"909\u00d78=\", \"170\u00d74=\"],\n  [\"677\u00d74=\", \"853\u00d79=\"],\n  [\"389\u00d72=\", \"474\u00d78=\"],\n  [\"188\u00d79=\", \"532\u00d74=\"],\n  [\"907\u00d73=\", \"285\u00d72=\"],\n  [\"124\u00d79=\", \"432\u00d79=\"],\n  [\"174\u00d76=\", \"947\u00d78=\"],\n  [\"844\u00d79=\", \"805\u00d79=\"],\n  [\"578\u00d77=\", \"467\u00d76=\"],\n  [\"525\u00d73=\", \"916\u00d74=\"],\n  [\"787\u00d74=\", \"284\u00d76=\"],\n  [\"503\u00d75=\", \"161\u00d79=\"],\n  [\"604\u00d77=\", \"483\u00d73=\"],\n  [\"895\u00d73=\", \"307\u00d75=\"],\n  [\"423\u00d74=\", \"193\u00d78=\"],\n  [\"275\u00d78=\", \"732\u00d73=\"],\n  [\"910\u00d73=\", \"402\u00d76=\"],\n  [\"677\u00d78=\", \"849\u00d76=\"],\n  [\"133\u00d76=\", \"523\u00d79=\"],\n  [\"616\u00d76=\", \"503\u00d73=\"],\n  [\"360\u00d75=\", \"240\u00d79=\"],\n  [\"440\u00d72=\", \"603\u00d75=\"],\n  [\"704\u00d77=\", \"276\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 multiplication problems to the\n# next day's values, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-15 Saturday\", \"2025-11-16 Sunday\"),\n    @(\"810\u00d77=\", \"718\u00d79=\"),\n    @(\"659\u00d79=\", \"633\u00d77=\"),\n    @(\"909\u00d78=\", \"170\u00d74=\"),\n    @(\"677\u00d74=\", \"853\u00d79=\"),\n    @(\"389\u00d72=\", \"474\u00d78=\"),\n    @(\"188\u00d79=\", \"532\u00d74=\"),\n    @(\"907\u00d73=\", \"285\u00d72=\"),\n    @(\"124\u00d79=\", \"432\u00d79=\"),\n    @(\"174\u00d76=\", \"947\u00d78=\"),\n    @(\"844\u00d79=\", \"805\u00d79=\"),\n    @(\"578\u00d77=\", \"467\u00d76=\"),\n    @(\"525\u00d73=\", \"916\u00d74=\"),\n    @(\"787\u00d74=\", \"284\u00d76=\"),\n    @(\"503\u00d75=\", \"161\u00d79=\"),\n    @(\"604\u00d77=\", \"483\u00d73=\"),\n    @(\"895\u00d73=\", \"307\u00d75=\"),\n    @(\"423\u00d74=\", \"193\u00d78=\"),\n    @(\"275\u00d78=\", \"732\u00d73=\"),\n    @(\"910\u00d73=\", \"402\u00d76=\"),\n    @(\"677\u00d78=\", \"849\u00d76=\"),\n    @(\"133\u00d76=\", \"523\u00d79=\"),\n    @(\"616\u00d76=\", \"503\u00d73=\"),\n    @(\"360\u00d75=\", \"240\u00d79=\"),\n    @(\"440\u00d72=\", \"603\u00d75=\"),\n    @(\"704\u00d77=\", \"276\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
